$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns D/E store plain text (inline strings) in the source
# data even when they look numeric (e.g. prices, percentages). When a
# value that parses as a plain number is assigned through .Value, Excel
# auto-converts the cell to a Number; to keep it text (matching the
# original authoring) we temporarily force a Text number format, then
# restore the cells original style so no visible formatting changes.

$ws.Range("D2").Value = "26.922.48"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "1.665.06"
$ws.Range("E3").Value = "  +0.57%  "
$origStyle_D5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.61"
$ws.Range("D5").Style = $origStyle_D5
$ws.Range("E5").Value = "  +0.25%  "
$origStyle_D6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.533"
$ws.Range("D6").Style = $origStyle_D6
$ws.Range("E6").Value = "  +4.78%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$origStyle_D8 = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0620"
$ws.Range("D8").Style = $origStyle_D8
$ws.Range("E8").Value = "  +0.88%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$origStyle_D9 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.249"
$ws.Range("D9").Style = $origStyle_D9
$ws.Range("E9").Value = "  +0.25%  "
$origStyle_D10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.29"
$ws.Range("D10").Style = $origStyle_D10
$ws.Range("E10").Value = "  +3.20%  "
$ws.Range("E11").Value = "  +3.89%  "
$ws.Range("D12").Value = "1.900.97"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").Value = "1.671.99"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("E15").Value = "  +0.83%  "
$origStyle_D16 = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.21"
$ws.Range("D16").Style = $origStyle_D16
$ws.Range("E16").Value = "  +2.04%  "
$ws.Range("D17").Value = "26.928.93"
$ws.Range("E17").Value = "  -0.44%  "
$origStyle_D18 = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.77"
$ws.Range("D18").Style = $origStyle_D18
$ws.Range("E18").Value = "  -1.88%  "
$origStyle_D19 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.99"
$ws.Range("D19").Style = $origStyle_D19
$ws.Range("E19").Value = "  +1.33%  "
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("E23").Value = "  -1.06%  "
$origStyle_D24 = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.13"
$ws.Range("D24").Style = $origStyle_D24
$ws.Range("E24").Value = "  -1.46%  "
$origStyle_D25 = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.10"
$ws.Range("D25").Style = $origStyle_D25
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("E27").Value = "  +1.40%  "
$origStyle_D28 = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.91"
$ws.Range("D28").Style = $origStyle_D28
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("E31").Value = "  -0.05%  "
$origStyle_D32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.36"
$ws.Range("D32").Style = $origStyle_D32
$ws.Range("E32").Value = "  +1.96%  "
$ws.Range("D33").Value = "1.458.65"
$ws.Range("E33").Value = "  -4.58%  "
$ws.Range("E34").Value = "  +3.24%  "
$ws.Range("E35").Value = "  +3.28%  "
$ws.Range("E36").Value = "  -0.16%  "
$origStyle_D37 = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.579"
$ws.Range("D37").Style = $origStyle_D37
$ws.Range("E37").Value = "  +0.30%  "
$origStyle_D38 = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.904"
$ws.Range("D38").Style = $origStyle_D38
$ws.Range("E38").Value = "  +1.82%  "
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("E40").Value = "  -3.25%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E42").Value = "  +0.49%  "
$origStyle_D43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.977"
$ws.Range("D43").Style = $origStyle_D43
$ws.Range("E43").Value = "  +6.12%  "
$origStyle_D44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "65.91"
$ws.Range("D44").Style = $origStyle_D44
$ws.Range("E44").Value = "  -0.65%  "
$ws.Range("D45").Value = "1.808.98"
$ws.Range("E45").Value = "  +0.74%  "
$ws.Range("E46").Value = "  +0.96%  "
$origStyle_D47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.44"
$ws.Range("D47").Style = $origStyle_D47
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("E50").Value = "  +3.96%  "
$ws.Range("E51").Value = "  +0.24%  "
